# Refresh cryptos list: updated Price (D) / Volume(1h) (E) values, and for
# rows 35-41 the coins were re-ranked, so Coin (B) / Link (C) / Price (D) /
# Volume(1h) (E) all shift to their new row. Index column A is left as-is.
#
# Numeric-looking Price strings (e.g. "225.36") are written with a leading
# apostrophe so Excel keeps them as text instead of auto-converting them to
# floating point numbers (matching the original inline-string cell type).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.681.99'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '1.810.18'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''225.36'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  +1.69%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '''40.80'
$ws.Range('E8').Value = '  +13.14%  '
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('E11').Value = '  +3.91%  '
$ws.Range('D12').Value = '2.071.69'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Value = '1.813.55'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '''10.94'
$ws.Range('E14').Value = '  -2.92%  '
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '34.678.30'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '''4.42'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').Value = '''67.99'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('D19').Value = '''241.91'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').Value = '0.0₃0771'
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').Value = '''11.13'
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '''4.11'
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('D25').Value = '''172.29'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('D26').Value = '''7.74'
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('D27').Value = '''17.49'
$ws.Range('E27').Value = '  +1.40%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '''3.79'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').Value = '''1.23'
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').Value = '''3.86'
$ws.Range('E33').Value = '  -1.51%  '
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''0.644'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = '''15.17'
$ws.Range('E36').Value = '  +13.75%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.313.59'
$ws.Range('E37').Value = '  -4.31%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.06'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('D39').Value = '''85.46'
$ws.Range('E39').Value = '  +4.52%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''2.37'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.0189'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('E42').Value = '  +6.66%  '
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').Value = '''0.945'
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('D46').Value = '''0.0520'
$ws.Range('E46').Value = '  +4.92%  '
$ws.Range('D47').Value = '1.970.27'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '''5.76'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '''101.27'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('E51').Value = '  +0.94%  '
